$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Generische Lebensmittel")

# Decrement column A values (rows 2..238) by 1 so the IDs start from 0
for ($r = 2; $r -le 238; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 - 1
}

# Update frozen pane / selection to point back at the top of the sheet
$ws.Range("A2").Select()
